$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '2026-02-27 07:18:37'
$ws.Range("N2").NumberFormat = "@"
$ws.Range("N2").Value = '0.5 °C 6:52 TU'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '2026-02-27 07:18:39'
$ws.Range("N3").NumberFormat = "@"
$ws.Range("N3").Value = '1.6 °C 6:51 TU'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '2026-02-27 07:18:42'
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = '100%'
$ws.Range("L4").NumberFormat = "@"
$ws.Range("L4").Value = '10.8 km/h - 315º 6:49 TU'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '2026-02-27 07:18:44'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '2026-02-27 07:18:46'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '2026-02-27 07:18:49'
$ws.Range("K7").NumberFormat = "@"
$ws.Range("K7").Value = '0.0 MJ/m2'
$ws.Range("N7").NumberFormat = "@"
$ws.Range("N7").Value = '8.3 °C 6:44 TU'
$ws.Range("O7").NumberFormat = "@"
$ws.Range("O7").Value = '9.8 °C'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '2026-02-27 07:18:51'
$ws.Range("L8").NumberFormat = "@"
$ws.Range("L8").Value = '14.8 km/h - 260º 6:36 TU'
$ws.Range("O8").NumberFormat = "@"
$ws.Range("O8").Value = '11.9 °C'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '2026-02-27 07:18:54'
$ws.Range("M9").NumberFormat = "@"
$ws.Range("M9").Value = '9.5 °C 6:56 TU'
$ws.Range("O9").NumberFormat = "@"
$ws.Range("O9").Value = '8.6 °C'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '2026-02-27 07:18:57'
$ws.Range("O10").NumberFormat = "@"
$ws.Range("O10").Value = '8.9 °C'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '2026-02-27 07:18:59'
$ws.Range("N11").NumberFormat = "@"
$ws.Range("N11").Value = '1.1 °C 6:44 TU'
$ws.Range("O11").NumberFormat = "@"
$ws.Range("O11").Value = '2.1 °C'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '2026-02-27 07:19:02'
$ws.Range("O12").NumberFormat = "@"
$ws.Range("O12").Value = '8.5 °C'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '2026-02-27 07:19:04'
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = '94%'
$ws.Range("K13").NumberFormat = "@"
$ws.Range("K13").Value = '0.0 MJ/m2'
$ws.Range("N13").NumberFormat = "@"
$ws.Range("N13").Value = '-3.9 °C 6:55 TU'
$ws.Range("O13").NumberFormat = "@"
$ws.Range("O13").Value = '-1.8 °C'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '2026-02-27 07:19:06'
$ws.Range("O14").NumberFormat = "@"
$ws.Range("O14").Value = '8.4 °C'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '2026-02-27 07:19:09'
$ws.Range("O15").NumberFormat = "@"
$ws.Range("O15").Value = '8.7 °C'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '2026-02-27 07:19:11'
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = '27%'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '2026-02-27 07:19:14'
$ws.Range("H17").NumberFormat = "@"
$ws.Range("H17").Value = '34%'
$ws.Range("K17").NumberFormat = "@"
$ws.Range("K17").Value = '0.0 MJ/m2'
$ws.Range("O17").NumberFormat = "@"
$ws.Range("O17").Value = '7.0 °C'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '2026-02-27 07:19:16'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '2026-02-27 07:19:19'
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = '82%'
$ws.Range("O19").NumberFormat = "@"
$ws.Range("O19").Value = '7.7 °C'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '2026-02-27 07:19:21'
$ws.Range("H20").NumberFormat = "@"
$ws.Range("H20").Value = '57%'
$ws.Range("K20").NumberFormat = "@"
$ws.Range("K20").Value = '0.0 MJ/m2'
$ws.Range("O20").NumberFormat = "@"
$ws.Range("O20").Value = '2.1 °C'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '2026-02-27 07:19:24'
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = '81%'
$ws.Range("J21").NumberFormat = "@"
$ws.Range("J21").Value = '1029.3 hPa'
$ws.Range("K21").NumberFormat = "@"
$ws.Range("K21").Value = '0.0 MJ/m2'
$ws.Range("N21").NumberFormat = "@"
$ws.Range("N21").Value = '1.5 °C 6:53 TU'
$ws.Range("O21").NumberFormat = "@"
$ws.Range("O21").Value = '3.3 °C'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '2026-02-27 07:19:26'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '2026-02-27 07:19:29'
$ws.Range("H23").NumberFormat = "@"
$ws.Range("H23").Value = '40%'
$ws.Range("O23").NumberFormat = "@"
$ws.Range("O23").Value = '2.8 °C'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '2026-02-27 07:19:31'
$ws.Range("J24").NumberFormat = "@"
$ws.Range("J24").Value = '1026.2 hPa'
$ws.Range("O24").NumberFormat = "@"
$ws.Range("O24").Value = '3.8 °C'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '2026-02-27 07:19:33'
$ws.Range("K25").NumberFormat = "@"
$ws.Range("K25").Value = '0.0 MJ/m2'
$ws.Range("O25").NumberFormat = "@"
$ws.Range("O25").Value = '4.6 °C'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '2026-02-27 07:19:36'
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = '2 cm'
$ws.Range("J26").NumberFormat = "@"
$ws.Range("J26").Value = '1024.4 hPa'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '2026-02-27 07:19:39'
$ws.Range("H27").NumberFormat = "@"
$ws.Range("H27").Value = '43%'
$ws.Range("M27").NumberFormat = "@"
$ws.Range("M27").Value = '4.9 °C 6:59 TU'
$ws.Range("O27").NumberFormat = "@"
$ws.Range("O27").Value = '4.2 °C'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '2026-02-27 07:19:41'
$ws.Range("L28").NumberFormat = "@"
$ws.Range("L28").Value = '10.1 km/h - 284º 6:59 TU'
$ws.Range("N28").NumberFormat = "@"
$ws.Range("N28").Value = '4.1 °C 6:41 TU'
$ws.Range("O28").NumberFormat = "@"
$ws.Range("O28").Value = '5.3 °C'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '2026-02-27 07:19:43'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '2026-02-27 07:19:46'
$ws.Range("O30").NumberFormat = "@"
$ws.Range("O30").Value = '9.7 °C'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '2026-02-27 07:19:48'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '2026-02-27 07:19:51'
$ws.Range("H32").NumberFormat = "@"
$ws.Range("H32").Value = '93%'
$ws.Range("N32").NumberFormat = "@"
$ws.Range("N32").Value = '-1.6 °C 6:43 TU'
$ws.Range("O32").NumberFormat = "@"
$ws.Range("O32").Value = '0.6 °C'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '2026-02-27 07:19:54'
$ws.Range("H33").NumberFormat = "@"
$ws.Range("H33").Value = '72%'
$ws.Range("N33").NumberFormat = "@"
$ws.Range("N33").Value = '0.3 °C 6:36 TU'
$ws.Range("O33").NumberFormat = "@"
$ws.Range("O33").Value = '2.0 °C'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '2026-02-27 07:19:56'
$ws.Range("L34").NumberFormat = "@"
$ws.Range("L34").Value = '21.2 km/h - 25º 6:39 TU'
$ws.Range("O34").NumberFormat = "@"
$ws.Range("O34").Value = '2.3 °C'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '2026-02-27 07:19:59'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '2026-02-27 07:20:01'
$ws.Range("M36").NumberFormat = "@"
$ws.Range("M36").Value = '10.6 °C 6:59 TU'
$ws.Range("O36").NumberFormat = "@"
$ws.Range("O36").Value = '9.5 °C'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '2026-02-27 07:20:04'
$ws.Range("N37").NumberFormat = "@"
$ws.Range("N37").Value = '0.0 °C 6:59 TU'
$ws.Range("O37").NumberFormat = "@"
$ws.Range("O37").Value = '2.4 °C'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '2026-02-27 07:20:07'
$ws.Range("O38").NumberFormat = "@"
$ws.Range("O38").Value = '7.4 °C'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '2026-02-27 07:20:09'
$ws.Range("K39").NumberFormat = "@"
$ws.Range("K39").Value = '0.0 MJ/m2'
$ws.Range("O39").NumberFormat = "@"
$ws.Range("O39").Value = '5.0 °C'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '2026-02-27 07:20:12'
$ws.Range("H40").NumberFormat = "@"
$ws.Range("H40").Value = '99%'
$ws.Range("N40").NumberFormat = "@"
$ws.Range("N40").Value = '0.4 °C 6:58 TU'
$ws.Range("O40").NumberFormat = "@"
$ws.Range("O40").Value = '1.7 °C'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '2026-02-27 07:20:14'
$ws.Range("J41").NumberFormat = "@"
$ws.Range("J41").Value = '1025.7 hPa'
$ws.Range("N41").NumberFormat = "@"
$ws.Range("N41").Value = '5.9 °C 6:31 TU'
$ws.Range("O41").NumberFormat = "@"
$ws.Range("O41").Value = '8.3 °C'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '2026-02-27 07:20:17'
$ws.Range("M42").NumberFormat = "@"
$ws.Range("M42").Value = '9.5 °C 6:58 TU'
$ws.Range("O42").NumberFormat = "@"
$ws.Range("O42").Value = '8.6 °C'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '2026-02-27 07:20:19'
$ws.Range("K43").NumberFormat = "@"
$ws.Range("K43").Value = '0.0 MJ/m2'
$ws.Range("N43").NumberFormat = "@"
$ws.Range("N43").Value = '1.6 °C 6:32 TU'
$ws.Range("O43").NumberFormat = "@"
$ws.Range("O43").Value = '3.6 °C'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '2026-02-27 07:20:22'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '2026-02-27 07:20:25'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '2026-02-27 07:20:27'
$ws.Range("O46").NumberFormat = "@"
$ws.Range("O46").Value = '6.2 °C'
